$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> new value for column F (dSF)
$updates = @{
    4  = 4
    8  = 2
    12 = 1
    16 = 3
    17 = -2
    27 = 2
    29 = 2
    33 = 0
    34 = 4
    35 = 4
    44 = -1
    46 = 0
    48 = 2
    56 = -4
    59 = -4
    65 = -2
    73 = 1
    82 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
